$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 186: 四方坪站充电量(kw), date 2025-12-02 (serial 45993)
$ws.Range("A186").Value = 45993
$ws.Range("B186").Value = "四方坪站充电量(kw)"
$ws.Range("C186").Value = 610.3180000000001
$ws.Range("D186").Value = 875.9030000000001
$ws.Range("E186").Value = 390.601
$ws.Range("F186").Value = 322.28
$ws.Range("G186").Value = 360.14300000000003
$ws.Range("H186").Value = 490.8990000000001
$ws.Range("I186").Value = 471.6479999999999
$ws.Range("J186").Value = 145.825
$ws.Range("K186").Value = 86.32
$ws.Range("L186").Value = 114.99
$ws.Range("M186").Value = 133.829
$ws.Range("N186").Value = 283.05
$ws.Range("O186").Value = 675.5830000000002
$ws.Range("P186").Value = 1297.7519999999997
$ws.Range("Q186").Value = 394.16999999999996
$ws.Range("R186").Value = 291.182
$ws.Range("S186").Value = 372.949
$ws.Range("T186").Value = 309.786
$ws.Range("U186").Value = 52.863
$ws.Range("V186").Value = 45.6
$ws.Range("W186").Value = 124.80000000000001
$ws.Range("X186").Value = 92.3
$ws.Range("Y186").Value = 32.38
$ws.Range("Z186").Value = 0

# Row 187: 高岭站充电量(kw), date 2025-12-02 (serial 45993)
$ws.Range("A187").Value = 45993
$ws.Range("B187").Value = "高岭站充电量(kw)"
$ws.Range("C187").Value = 251.09199999999998
$ws.Range("D187").Value = 259.71
$ws.Range("E187").Value = 69.98100000000001
$ws.Range("F187").Value = 213.125
$ws.Range("G187").Value = 45.294
$ws.Range("H187").Value = 156.339
$ws.Range("I187").Value = 318.063
$ws.Range("J187").Value = 132.839
$ws.Range("K187").Value = 251.82700000000003
$ws.Range("L187").Value = 327.147
$ws.Range("M187").Value = 339.10699999999997
$ws.Range("N187").Value = 381.40000000000003
$ws.Range("O187").Value = 707.0020000000001
$ws.Range("P187").Value = 661.8689999999999
$ws.Range("Q187").Value = 341.5469999999999
$ws.Range("R187").Value = 537.978
$ws.Range("S187").Value = 260.56600000000003
$ws.Range("T187").Value = 111.493
$ws.Range("U187").Value = 13.422
$ws.Range("V187").Value = 81.56200000000001
$ws.Range("W187").Value = 44.995999999999995
$ws.Range("X187").Value = 1.945
$ws.Range("Y187").Value = 17.445
$ws.Range("Z187").Value = 88.731

# Update selection to reflect new active cell as in the diff
$ws.Range("D191").Select()
